$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0
$ws.Cells.Item(2, 8).Value = 0.09938754597967438
$ws.Range("C2:G2").ClearContents()

$ws.Cells.Item(3, 2).Value = 0.09544247632357654
$ws.Cells.Item(3, 8).Value = 0.1948300223032509
$ws.Range("C3:G3").ClearContents()

$ws.Cells.Item(4, 2).Value = 0.07433551145063542
$ws.Cells.Item(4, 8).Value = 0.1737230574303098
$ws.Range("C4:G4").ClearContents()

$ws.Cells.Item(5, 2).Value = 0.04606018388745416
$ws.Cells.Item(5, 8).Value = 0.1454477298671285
$ws.Range("C5:G5").ClearContents()

$ws.Cells.Item(6, 2).Value = 0.06073887638649635
$ws.Cells.Item(6, 8).Value = 0.1601264223661707
$ws.Range("C6:G6").ClearContents()

$ws.Cells.Item(7, 2).Value = 0.05558531632287125
$ws.Cells.Item(7, 8).Value = 0.1549728623025456
$ws.Range("C7:G7").ClearContents()

$ws.Cells.Item(8, 2).Value = 0.03768428929135112
$ws.Cells.Item(8, 8).Value = 0.1370718352710255
$ws.Cells.Item(8, 3).Value = 0.002056094966887085
$ws.Cells.Item(8, 4).Value = 4.610602075245989
$ws.Cells.Item(8, 5).Value = 0.0145225967495337
$ws.Cells.Item(8, 6).Value = 0.03365153501609233
$ws.Cells.Item(8, 7).Value = 0.04171704356661073

$ws.Cells.Item(9, 2).Value = 0.03702959467894879
$ws.Cells.Item(9, 8).Value = 0.1364171406586231
$ws.Range("C9:G9").ClearContents()

$ws.Cells.Item(10, 2).Value = 0.03750736443914262
$ws.Cells.Item(10, 8).Value = 0.136894910418817
$ws.Cells.Item(10, 3).Value = 0.002441190054921404
$ws.Cells.Item(10, 4).Value = 4.55487105918922
$ws.Cells.Item(10, 5).Value = 0.00892760710090874
$ws.Cells.Item(10, 6).Value = 0.03271774425642111
$ws.Cells.Item(10, 7).Value = 0.04229698462186376

$ws.Cells.Item(11, 2).Value = 0.0299333419793639
$ws.Cells.Item(11, 8).Value = 0.1293208879590383
$ws.Range("C11:G11").ClearContents()

$ws.Cells.Item(12, 2).Value = 0.05235207103284203
$ws.Cells.Item(12, 8).Value = 0.1517396170125164
$ws.Range("C12:G12").ClearContents()

$ws.Cells.Item(13, 2).Value = 0.06459953424568558
$ws.Cells.Item(13, 8).Value = 0.16398708022536
$ws.Range("C13:G13").ClearContents()

$ws.Cells.Item(14, 2).Value = 0.07541470999723525
$ws.Cells.Item(14, 8).Value = 0.1748022559769096
$ws.Range("C14:G14").ClearContents()

$ws.Cells.Item(15, 2).Value = 0.08077247944632321
$ws.Cells.Item(15, 8).Value = 0.1801600254259976
$ws.Range("C15:G15").ClearContents()

$ws.Cells.Item(16, 2).Value = 0.08348292184851812
$ws.Cells.Item(16, 8).Value = 0.1828704678281925
$ws.Range("C16:G16").ClearContents()

$ws.Cells.Item(17, 2).Value = 0.08625249812131089
$ws.Cells.Item(17, 8).Value = 0.1856400441009853
$ws.Range("C17:G17").ClearContents()

$ws.Cells.Item(18, 2).Value = -0.09938754597967438
$ws.Cells.Item(18, 8).Value = 0
$ws.Cells.Item(18, 3).Value = 0.008653689853673477
$ws.Cells.Item(18, 4).Value = -18.86719455148055
$ws.Cells.Item(18, 5).Value = 0.0246845553515376
$ws.Cells.Item(18, 6).Value = -0.1164063896228987
$ws.Cells.Item(18, 7).Value = -0.08236870233644973

$ws.Cells.Item(19, 2).Value = 0.09066174891765244
$ws.Cells.Item(19, 8).Value = 0.1900492948973268
$ws.Range("C19:G19").ClearContents()

$ws.Cells.Item(20, 2).Value = 0.09232087629706269
$ws.Cells.Item(20, 8).Value = 0.1917084222767371
$ws.Range("C20:G20").ClearContents()

$ws.Cells.Item(21, 2).Value = 0.09727129600694094
$ws.Cells.Item(21, 8).Value = 0.1966588419866153
$ws.Range("C21:G21").ClearContents()

$ws.Cells.Item(22, 2).Value = 0.1000947009958603
$ws.Cells.Item(22, 8).Value = 0.1994822469755347
$ws.Range("C22:G22").ClearContents()

$ws.Cells.Item(23, 2).Value = 0.103999561416143
$ws.Cells.Item(23, 8).Value = 0.2033871073958173
$ws.Cells.Item(23, 3).Value = 0.007785622198110591
$ws.Cells.Item(23, 4).Value = -323436960025.2491
$ws.Cells.Item(23, 5).Value = 0.04744812292787406
$ws.Cells.Item(23, 6).Value = 0.08869425653373132
$ws.Cells.Item(23, 7).Value = 0.1193048662985546

$ws.Cells.Item(24, 2).Value = 0.1068487404761665
$ws.Cells.Item(24, 8).Value = 0.2062362864558409
$ws.Range("C24:G24").ClearContents()

$ws.Cells.Item(25, 2).Value = 0.1091724292302566
$ws.Cells.Item(25, 8).Value = 0.2085599752099309
$ws.Cells.Item(25, 3).Value = 0.007745187477226081
$ws.Cells.Item(25, 4).Value = 25.64992793640603
$ws.Cells.Item(25, 5).Value = 0.04992955927607987
$ws.Cells.Item(25, 6).Value = 0.09394820866560222
$ws.Cells.Item(25, 7).Value = 0.1243966497949112

$ws.Cells.Item(26, 2).Value = 0.1093807751548355
$ws.Cells.Item(26, 8).Value = 0.2087683211345099
$ws.Cells.Item(26, 3).Value = 0.007758347567464143
$ws.Cells.Item(26, 4).Value = -111759430900.8102
$ws.Cells.Item(26, 5).Value = 0.05997937961130436
$ws.Cells.Item(26, 6).Value = 0.09413546175669202
$ws.Cells.Item(26, 7).Value = 0.1246260885529791

$ws.Cells.Item(27, 2).Value = 0.1129574129246157
$ws.Cells.Item(27, 8).Value = 0.2123449589042901
$ws.Range("C27:G27").ClearContents()

$ws.Cells.Item(28, 2).Value = 0.1076206129339424
$ws.Cells.Item(28, 8).Value = 0.2070081589136168
$ws.Cells.Item(28, 3).Value = 0.007314701122807652
$ws.Cells.Item(28, 4).Value = 22.52232618684381
$ws.Cells.Item(28, 5).Value = 0.08773815096072336
$ws.Cells.Item(28, 6).Value = 0.09325684945791586
$ws.Cells.Item(28, 7).Value = 0.1219843764099697

$ws.Cells.Item(29, 2).Value = 0.04175810814375928
$ws.Cells.Item(29, 8).Value = 0.1411456541234337
$ws.Cells.Item(29, 3).Value = 0.002461833108853095
$ws.Cells.Item(29, 4).Value = 4.917267414069542
$ws.Cells.Item(29, 5).Value = 0.0105151434110451
$ws.Cells.Item(29, 6).Value = 0.03692382100225468
$ws.Cells.Item(29, 7).Value = 0.04659239528526291
